$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "Подарил" (D) / "Получил" (E) columns from "нет" to "да"
# for the specific rows that changed in the source table.
$ws.Range("D10").Value = "да"
$ws.Range("E10").Value = "да"

$ws.Range("E16").Value = "да"

$ws.Range("D46").Value = "да"

$ws.Range("D60").Value = "да"
$ws.Range("E60").Value = "да"

$ws.Range("D61").Value = "да"
$ws.Range("E61").Value = "да"

$ws.Range("D74").Value = "да"
